$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hashcode (column B) values for the rows whose code (column A)
# matches the entries below. Each assignment is guarded by a check that
# column A still holds the expected code, so the right row is touched
# even if something upstream shifted the data.

if ($ws.Cells.Item(9, 1).Text -eq "05-050305TC") {
    $ws.Cells.Item(9, 2).Value = "63ca7202aa753dbe027119bfeda10078"
}
if ($ws.Cells.Item(17, 1).Text -eq "05-050305TP") {
    $ws.Cells.Item(17, 2).Value = "40d14c4d3bdea334539e0f9099d74ead"
}
if ($ws.Cells.Item(29, 1).Text -eq "05-050302A") {
    $ws.Cells.Item(29, 2).Value = "bbe541ad11657fcb5ec5ac6d37dc99f2"
}
if ($ws.Cells.Item(34, 1).Text -eq "05-050316TP") {
    $ws.Cells.Item(34, 2).Value = "661cccda4f7697e3d266ad7cde7e4415"
}
if ($ws.Cells.Item(126, 1).Text -eq "05-050309A") {
    $ws.Cells.Item(126, 2).Value = "f26ece442805a7105448d592e0c137da"
}
if ($ws.Cells.Item(136, 1).Text -eq "05-050312TC") {
    $ws.Cells.Item(136, 2).Value = "f6eae73292a6e03f62c8b30f714a6710"
}
if ($ws.Cells.Item(163, 1).Text -eq "05-050308A") {
    $ws.Cells.Item(163, 2).Value = "05be76d980c561e85bc5d139ef0796ea"
}
if ($ws.Cells.Item(176, 1).Text -eq "05-050303TP") {
    $ws.Cells.Item(176, 2).Value = "08a7efbd4a3864f167faa5eea348ecd2"
}
if ($ws.Cells.Item(181, 1).Text -eq "05-050303TC") {
    $ws.Cells.Item(181, 2).Value = "9aec6a43074d9dba65c7daab5f0f307b"
}
if ($ws.Cells.Item(184, 1).Text -eq "05-050305A") {
    $ws.Cells.Item(184, 2).Value = "ac53f209c4d8c0bcc8bde7656604161e"
}
if ($ws.Cells.Item(201, 1).Text -eq "05-050306A") {
    $ws.Cells.Item(201, 2).Value = "ae31a06e44f0140f50906977a258e83b"
}
if ($ws.Cells.Item(214, 1).Text -eq "05-050303A") {
    $ws.Cells.Item(214, 2).Value = "be4fc743f04690faa9a8135ed11d0c21"
}
if ($ws.Cells.Item(230, 1).Text -eq "05-050304A") {
    $ws.Cells.Item(230, 2).Value = "d63f1cba49810b71c8dee8a6872a13a6"
}
if ($ws.Cells.Item(248, 1).Text -eq "05-050003TC") {
    $ws.Cells.Item(248, 2).Value = "d49dab1ad20a4ad7828b85f6324ae86c"
}
if ($ws.Cells.Item(278, 1).Text -eq "01-080101-010112TM") {
    $ws.Cells.Item(278, 2).Value = "8c6e2b75376b8490b816902250befb49"
}
if ($ws.Cells.Item(282, 1).Text -eq "05-050003TP") {
    $ws.Cells.Item(282, 2).Value = "48ef184f805a5536e4de511cdc3e7ee4"
}
if ($ws.Cells.Item(299, 1).Text -eq "05-050310TC") {
    $ws.Cells.Item(299, 2).Value = "4d2d3338ed5bddb3bf594da95aa52b59"
}
if ($ws.Cells.Item(308, 1).Text -eq "05-050310TP") {
    $ws.Cells.Item(308, 2).Value = "28ff6935881dfd4de4ae62f37220508d"
}
if ($ws.Cells.Item(505, 1).Text -eq "05-050208TC") {
    $ws.Cells.Item(505, 2).Value = "12beed6b9d0c50af8787dfa8a664a090"
}
if ($ws.Cells.Item(513, 1).Text -eq "05-050202A") {
    $ws.Cells.Item(513, 2).Value = "fa6312013da5e9a41c9724b330b457e9"
}
if ($ws.Cells.Item(514, 1).Text -eq "05-050311A") {
    $ws.Cells.Item(514, 2).Value = "e6486037e88db1099dc96a3609e9c79d"
}
if ($ws.Cells.Item(520, 1).Text -eq "05-050306TP") {
    $ws.Cells.Item(520, 2).Value = "5e9804d8cbe33be6afba717160debf6d"
}
if ($ws.Cells.Item(528, 1).Text -eq "05-050317TC") {
    $ws.Cells.Item(528, 2).Value = "6edbd2c0738d73814d65d512c12c4639"
}
if ($ws.Cells.Item(539, 1).Text -eq "05-050317TP") {
    $ws.Cells.Item(539, 2).Value = "487b086f90b836403abd7d919bb7235c"
}
if ($ws.Cells.Item(563, 1).Text -eq "05-050201A") {
    $ws.Cells.Item(563, 2).Value = "348d0d3e4f4df3e7482aa7c6f230e00f"
}
if ($ws.Cells.Item(566, 1).Text -eq "05-050310A") {
    $ws.Cells.Item(566, 2).Value = "12a358635a841bf54a69664e6b694f9c"
}
if ($ws.Cells.Item(579, 1).Text -eq "05-050308TC") {
    $ws.Cells.Item(579, 2).Value = "73c1d342a327e32f561c83e276c591d3"
}
if ($ws.Cells.Item(582, 1).Text -eq "05-050004A") {
    $ws.Cells.Item(582, 2).Value = "d702cc955674adb5daa772a3e4032392"
}
if ($ws.Cells.Item(588, 1).Text -eq "05-050308TP") {
    $ws.Cells.Item(588, 2).Value = "86d6a105e7354519b7e65ae96e692316"
}
if ($ws.Cells.Item(600, 1).Text -eq "05-050005A") {
    $ws.Cells.Item(600, 2).Value = "30b32785d064080f176f41543296c20c"
}
if ($ws.Cells.Item(645, 1).Text -eq "05-050302TP") {
    $ws.Cells.Item(645, 2).Value = "2490f1df30605f9fee490ce3c88285be"
}
if ($ws.Cells.Item(682, 1).Text -eq "05-050317A") {
    $ws.Cells.Item(682, 2).Value = "9f3ab27c1fbcc29b3dc5a75d2ca0ece2"
}
if ($ws.Cells.Item(716, 1).Text -eq "05-050304TC") {
    $ws.Cells.Item(716, 2).Value = "1c722e64a9f30fc6968bdcb5f35b55e9"
}
if ($ws.Cells.Item(731, 1).Text -eq "05-050304TP") {
    $ws.Cells.Item(731, 2).Value = "ac14abdc33f59e962e2b9d00792a0c03"
}
if ($ws.Cells.Item(742, 1).Text -eq "05-050315TC") {
    $ws.Cells.Item(742, 2).Value = "5b35ff0db4bb7fa705668d76f0679e4e"
}
if ($ws.Cells.Item(745, 1).Text -eq "05-050316A") {
    $ws.Cells.Item(745, 2).Value = "455e9300c164959ef5eac100ff083875"
}
if ($ws.Cells.Item(772, 1).Text -eq "05-050004TC") {
    $ws.Cells.Item(772, 2).Value = "3b050dca3539ed9b0f1a03a757100a28"
}
if ($ws.Cells.Item(778, 1).Text -eq "05-050004TP") {
    $ws.Cells.Item(778, 2).Value = "70534bb0da705019895dfef9f87114ef"
}
if ($ws.Cells.Item(833, 1).Text -eq "05-050104TC") {
    $ws.Cells.Item(833, 2).Value = "dd61380f9e32c3c52edac4f3ab73c6af"
}
if ($ws.Cells.Item(836, 1).Text -eq "05-050202TP") {
    $ws.Cells.Item(836, 2).Value = "f98dad657e0ab667ef0d4f13e0ea9c43"
}
if ($ws.Cells.Item(842, 1).Text -eq "05-050104TM") {
    $ws.Cells.Item(842, 2).Value = "c27d43645588174d7ecedd33bda0fe5a"
}
if ($ws.Cells.Item(853, 1).Text -eq "05-050311TP") {
    $ws.Cells.Item(853, 2).Value = "331d72bc9756f8ec510f4366059629f5"
}
if ($ws.Cells.Item(872, 1).Text -eq "05-050309TC") {
    $ws.Cells.Item(872, 2).Value = "b5024786e229ecca267ebafc25c40b5c"
}
if ($ws.Cells.Item(880, 1).Text -eq "05-050309TP") {
    $ws.Cells.Item(880, 2).Value = "4f8d110331f465306f3197c522fea58c"
}
if ($ws.Cells.Item(887, 1).Text -eq "05-050003A") {
    $ws.Cells.Item(887, 2).Value = "fd1c7181a00eaa25edac3b8b203ad812"
}
if ($ws.Cells.Item(923, 1).Text -eq "05-050001A") {
    $ws.Cells.Item(923, 2).Value = "74a9676d72e412b6264b0fe7119d182e"
}
if ($ws.Cells.Item(947, 1).Text -eq "05-050002TP") {
    $ws.Cells.Item(947, 2).Value = "10785968df9b796fedd30771977d49e9"
}
